$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row for the latest collection date (2021-06-08, Tue)
#    above the existing most-recent-week row (row 6), pushing all the
#    historical rows down by one.
$ws.Rows.Item(6).Insert()

# 2) Populate the new row 6 with this week's figures.
$ws.Range("A6").Value = 44355
$ws.Range("B6").Value = "(火)"
$ws.Range("C6").Formula = "=SUM(D6:G6)"
$ws.Range("D6").Value = 71689
$ws.Range("E6").Value = 990
$ws.Range("F6").Value = 89243
$ws.Range("G6").Value = 0

# 3) Match the formatting used by the other weekly data rows (copy down
#    from the row directly below, which now holds the former row 6).
$ws.Range("A7:I7").Copy()
$ws.Range("A6:I6").PasteSpecial(-4122)
$ws.Range("A6").NumberFormat = "mm-dd-yy"

# 4) Roll the new week's counts into the running cumulative totals row.
$ws.Range("D5").Value = 5021299
$ws.Range("E5").Value = 6312
$ws.Range("F5").Value = 3628328

# 5) Update the "as of" footnote date shown at G2.
$ws.Range("G2").Value = "（6月8日時点）"

# 6) Extend the print area by the one additional row.
$ws.PageSetup.PrintArea = "A1:G49"
